$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update recipe header labels
$ws.Range("B1").Value = "Recipe Rice and Beef"
$ws.Range("C1").Value = "Recipe Rice and chicken"
$ws.Range("D1").Value = "Recipe Rice and Soy tofu"

# Swap the values in columns C and D for data rows 2-15
for ($r = 2; $r -le 15; $r++) {
    $cVal = $ws.Cells.Item($r, 3).Value2
    $dVal = $ws.Cells.Item($r, 4).Value2
    $ws.Cells.Item($r, 3).Value2 = $dVal
    $ws.Cells.Item($r, 4).Value2 = $cVal
}
